$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2024, "Hipertensão", 37068, 18044255.16, 486.79, 4.22),
    @(2014, "Doenças preveníveis por imunização e condições sensíveis", 17733, 26005097.38, 1466.48, 15.83),
    @(2014, "Doenças relacionadas ao Pré-Natal e Parto", 168, 261426.06, 1556.11, 8.34),
    @(2024, "Diabetes melitus", 136172, 151124585.96, 1109.81, 6.78),
    @(2024, "Úlcera gastrointestinal", 50844, 60521054.37, 1190.33, 5.6),
    @(2024, "Pneumonias bacterianas", 168655, 243172884.92, 1441.84, 6.83),
    @(2014, "Doenças Cerebrovasculares", 141909, 172055725.92, 1212.44, 7.72),
    @(2024, "Doenças Cerebrovasculares", 197963, 365759200.23, 1847.61, 7.29),
    @(2014, "Diabetes melitus", 139819, 89667747.68000001, 641.3099999999999, 6),
    @(2014, "Pneumonias bacterianas", 149040, 152952601.45, 1026.25, 6.27),
    @(2014, "Úlcera gastrointestinal", 26980, 20666487.39, 765.99, 5.38),
    @(2024, "Infecções de ouvido, nariz e garganta", 2239, 670022.02, 299.25, 2.89),
    @(2014, "Infecção no Rim e Trato Urinário", 186646, 63451630.38, 339.96, 4.47),
    @(2024, "Infecção da pele e tecido subcutâneo", 61046, 52273362.37, 856.29, 7.34),
    @(2024, "Doencas pulmonares", 22406, 20922481.99, 933.79, 4.71),
    @(2014, "Asma", 12292, 6598149.13, 536.78, 2.93),
    @(2014, "Doença Inflamatória órgãos pélvicos femininos", 309, 72067.13, 233.23, 1.7),
    @(2024, "Insuficiência Cardíaca", 26020, 71171979.81, 2735.28, 8.26),
    @(2014, "Gastroenterites Infecciosas e complicações", 204186, 67409507.34999999, 330.14, 3.11),
    @(2014, "Deficiências Nutricionais", 37174, 24500807.12, 659.08, 7.25),
    @(2024, "Infecção no Rim e Trato Urinário", 218241, 126877175.19, 581.36, 5.67),
    @(2024, "Gastroenterites Infecciosas e complicações", 126223, 51377873.76, 407.04, 3.33),
    @(2024, "Deficiências Nutricionais", 19316, 24227126.98, 1254.25, 8.67),
    @(2014, "Insuficiência Cardíaca", 13725, 24408951.27, 1778.43, 7.02),
    @(2024, "Eplepsias", 1485, 1473095.85, 991.98, 5.62),
    @(2024, "Doenças preveníveis por imunização e condições sensíveis", 23839, 46300850.19, 1942.23, 13.73),
    @(2014, "Hipertensão", 77945, 26713533.69, 342.72, 4.34),
    @(2024, "Anemia", 1121, 446850.07, 398.62, 3.99),
    @(2024, "Doenças relacionadas ao Pré-Natal e Parto", 1827, 1143375.44, 625.8200000000001, 6.28),
    @(2024, "Asma", 6492, 4849799.21, 747.04, 3.43),
    @(2014, "Doencas pulmonares", 18493, 12387346.96, 669.84, 4.09),
    @(2014, "Infecção da pele e tecido subcutâneo", 41557, 26165231.64, 629.62, 6.52),
    @(2024, "Doença Inflamatória órgãos pélvicos femininos", 1081, 480510.46, 444.51, 2.45),
    @(2024, "Angina", 1132, 5301099.39, 4682.95, 4.87),
    @(2014, "Infecções de ouvido, nariz e garganta", 84, 5391.6, 64.19, 1.52)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
    $ws.Cells.Item($r, 4).Value = $data[$i][3]
    $ws.Cells.Item($r, 5).Value = $data[$i][4]
    $ws.Cells.Item($r, 6).Value = $data[$i][5]
}

Write-Host "Updated rows 2-36"
